# Illinois_Converted workbook update:
#  - Weights row (row 7): travel_limit columns X/Y (0.5 -> 0), total weight AC7 (13 -> 12)
#  - Recompute the weighted-average column AC for every data row (9..221) using the
#    new weights / new total so the ratios stay consistent with the updated weights
#  - Append 12 new daily policy rows (222..233) for 9/30/2020 .. 10/11/2020, carrying
#    forward the same category flags as the last existing row (221)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Converted Data")

$firstDataCol = 2   # column B
$lastDataCol  = 28  # column AB
$weightCol    = 29  # column AC

# --- Update the weights row (row 7): travel_limit (X7/Y7) drop to 0 ---
$ws.Cells.Item(7, 24).Value = 0   # X7
$ws.Cells.Item(7, 25).Value = 0   # Y7

# Total weight (AC7) is stored as a clean round number (13 -> 12), not a live SUM
# formula (the underlying weights don't add up to an exact integer due to repeating
# decimals like 0.1666666667), so set it explicitly to match.
$totalWeight = 12
$ws.Cells.Item(7, $weightCol).Value = $totalWeight

# Cache the (now-updated) weights row for reuse below
$weights = @{}
for ($c = $firstDataCol; $c -le $lastDataCol; $c++) {
    $weights[$c] = $ws.Cells.Item(7, $c).Value2
}

# --- Recompute AC (weighted average) for every existing data row 9..221 ---
for ($r = 9; $r -le 221; $r++) {
    $sum = 0
    for ($c = $firstDataCol; $c -le $lastDataCol; $c++) {
        $cellVal = $ws.Cells.Item($r, $c).Value2
        if ($cellVal -ne $null) {
            $sum = $sum + ($cellVal * $weights[$c])
        }
    }
    $ws.Cells.Item($r, $weightCol).Value = $sum / $totalWeight
}

# --- Append 12 new rows (222..233) for 9/30/2020 .. 10/11/2020 ---
# Same per-category flags as row 221 (columns B..AB), style "bold/border" on column A
$newDates = @("9/30/2020","10/1/2020","10/2/2020","10/3/2020","10/4/2020","10/5/2020","10/6/2020","10/7/2020","10/8/2020","10/9/2020","10/10/2020","10/11/2020")

$templateRow = 221
$startRow = 222

for ($i = 0; $i -lt $newDates.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newDates[$i]
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 1)).Style = $ws.Range($ws.Cells.Item($templateRow, 1), $ws.Cells.Item($templateRow, 1)).Style
    for ($c = $firstDataCol; $c -le $lastDataCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($templateRow, $c).Value2
    }
    $ws.Cells.Item($r, $weightCol).Value = 0.25
}
